$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 5.57196
$ws.Cells.Item(2, 8).Value = 16.71588
$ws.Cells.Item(2, 9).Value = 0.2041274261050298
$ws.Cells.Item(2, 10).Value = 0.2041274261050299
$ws.Cells.Item(2, 13).Value = 15.090721
$ws.Cells.Item(2, 14).Value = 45.272163
$ws.Cells.Item(2, 15).Value = 0.169971412714946
$ws.Cells.Item(2, 16).Value = 0.169971412714946
$ws.Cells.Item(2, 17).Value = 84.08489378316
$ws.Cells.Item(2, 18).Value = 756.7640440484399
$ws.Cells.Item(2, 19).Value = 0.03469582698893767
$ws.Cells.Item(2, 20).Value = 0.03469582698893767

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 5.57196
$ws.Cells.Item(3, 8).Value = 16.71588
$ws.Cells.Item(3, 9).Value = 0.2041274261050298
$ws.Cells.Item(3, 10).Value = 0.2041274261050299
$ws.Cells.Item(3, 14).Value = 60.45961299999999
$ws.Cells.Item(3, 15).Value = 0.226991713071207
$ws.Cells.Item(3, 16).Value = 0.226991713071207
$ws.Cells.Item(3, 17).Value = 112.29284841716
$ws.Cells.Item(3, 18).Value = 1010.63563575444
$ws.Cells.Item(3, 19).Value = 0.04633523413639695
$ws.Cells.Item(3, 20).Value = 0.04633523413639695

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 5.57196
$ws.Cells.Item(4, 8).Value = 16.71588
$ws.Cells.Item(4, 9).Value = 0.2041274261050298
$ws.Cells.Item(4, 10).Value = 0.2041274261050299
$ws.Cells.Item(4, 13).Value = 18.66868666666667
$ws.Cells.Item(4, 14).Value = 56.00606
$ws.Cells.Item(4, 15).Value = 0.2102711359030499
$ws.Cells.Item(4, 16).Value = 0.2102711359030499
$ws.Cells.Item(4, 17).Value = 104.0211753592
$ws.Cells.Item(4, 18).Value = 936.1905782327999
$ws.Cells.Item(4, 19).Value = 0.04292210575607051
$ws.Cells.Item(4, 20).Value = 0.04292210575607052

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 5.57196
$ws.Cells.Item(5, 8).Value = 16.71588
$ws.Cells.Item(5, 9).Value = 0.2041274261050298
$ws.Cells.Item(5, 10).Value = 0.2041274261050299
$ws.Cells.Item(5, 13).Value = 5.641943333333334
$ws.Cells.Item(5, 14).Value = 16.92583
$ws.Cells.Item(5, 15).Value = 0.06354693581733691
$ws.Cells.Item(5, 16).Value = 0.06354693581733692
$ws.Cells.Item(5, 17).Value = 31.4366825756
$ws.Cells.Item(5, 18).Value = 282.9301431804
$ws.Cells.Item(5, 19).Value = 0.01297167244525451
$ws.Cells.Item(5, 20).Value = 0.01297167244525452

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 5.57196
$ws.Cells.Item(6, 8).Value = 16.71588
$ws.Cells.Item(6, 9).Value = 0.2041274261050298
$ws.Cells.Item(6, 10).Value = 0.2041274261050299
$ws.Cells.Item(6, 13).Value = 15.42507366666667
$ws.Cells.Item(6, 14).Value = 46.275221
$ws.Cells.Item(6, 15).Value = 0.1737373292074942
$ws.Cells.Item(6, 16).Value = 0.1737373292074942
$ws.Cells.Item(6, 17).Value = 85.94789346772001
$ws.Cells.Item(6, 18).Value = 773.53104120948
$ws.Cells.Item(6, 19).Value = 0.03546455382948801
$ws.Cells.Item(6, 20).Value = 0.03546455382948802

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 5.57196
$ws.Cells.Item(7, 8).Value = 16.71588
$ws.Cells.Item(7, 9).Value = 0.2041274261050298
$ws.Cells.Item(7, 10).Value = 0.2041274261050299
$ws.Cells.Item(7, 13).Value = 13.804248
$ws.Cells.Item(7, 14).Value = 41.412744
$ws.Cells.Item(7, 15).Value = 0.1554814732859661
$ws.Cells.Item(7, 16).Value = 0.1554814732859661
$ws.Cells.Item(7, 17).Value = 76.91671768607999
$ws.Cells.Item(7, 18).Value = 692.2504591747199
$ws.Cells.Item(7, 19).Value = 0.03173803294888222
$ws.Cells.Item(7, 20).Value = 0.03173803294888223

$ws.Cells.Item(8, 7).Value = 20.56891366666667
$ws.Cells.Item(8, 8).Value = 61.70674100000001
$ws.Cells.Item(8, 9).Value = 0.7535372480335895
$ws.Cells.Item(8, 10).Value = 0.7535372480335896
$ws.Cells.Item(8, 13).Value = 15.090721
$ws.Cells.Item(8, 14).Value = 45.272163
$ws.Cells.Item(8, 15).Value = 0.169971412714946
$ws.Cells.Item(8, 16).Value = 0.169971412714946
$ws.Cells.Item(8, 17).Value = 310.3997374167537
$ws.Cells.Item(8, 18).Value = 2793.597636750783
$ws.Cells.Item(8, 19).Value = 0.1280797905816019
$ws.Cells.Item(8, 20).Value = 0.1280797905816019

$ws.Cells.Item(9, 7).Value = 20.56891366666667
$ws.Cells.Item(9, 8).Value = 61.70674100000001
$ws.Cells.Item(9, 9).Value = 0.7535372480335895
$ws.Cells.Item(9, 10).Value = 0.7535372480335896
$ws.Cells.Item(9, 14).Value = 60.45961299999999
$ws.Cells.Item(9, 15).Value = 0.226991713071207
$ws.Cells.Item(9, 16).Value = 0.226991713071207
$ws.Cells.Item(9, 17).Value = 414.5295200390259
$ws.Cells.Item(9, 18).Value = 3730.765680351233
$ws.Cells.Item(9, 19).Value = 0.1710467107941075
$ws.Cells.Item(9, 20).Value = 0.1710467107941075

$ws.Cells.Item(10, 7).Value = 20.56891366666667
$ws.Cells.Item(10, 8).Value = 61.70674100000001
$ws.Cells.Item(10, 9).Value = 0.7535372480335895
$ws.Cells.Item(10, 10).Value = 0.7535372480335896
$ws.Cells.Item(10, 13).Value = 18.66868666666667
$ws.Cells.Item(10, 14).Value = 56.00606
$ws.Cells.Item(10, 15).Value = 0.2102711359030499
$ws.Cells.Item(10, 16).Value = 0.2102711359030499
$ws.Cells.Item(10, 17).Value = 383.9946043167178
$ws.Cells.Item(10, 18).Value = 3455.95143885046
$ws.Cells.Item(10, 19).Value = 0.1584471330892811
$ws.Cells.Item(10, 20).Value = 0.1584471330892812

$ws.Cells.Item(11, 7).Value = 20.56891366666667
$ws.Cells.Item(11, 8).Value = 61.70674100000001
$ws.Cells.Item(11, 9).Value = 0.7535372480335895
$ws.Cells.Item(11, 10).Value = 0.7535372480335896
$ws.Cells.Item(11, 13).Value = 5.641943333333334
$ws.Cells.Item(11, 14).Value = 16.92583
$ws.Cells.Item(11, 15).Value = 0.06354693581733691
$ws.Cells.Item(11, 16).Value = 0.06354693581733692
$ws.Cells.Item(11, 17).Value = 116.0486453355589
$ws.Cells.Item(11, 18).Value = 1044.43780802003
$ws.Cells.Item(11, 19).Value = 0.04788498313676319
$ws.Cells.Item(11, 20).Value = 0.0478849831367632

$ws.Cells.Item(12, 7).Value = 20.56891366666667
$ws.Cells.Item(12, 8).Value = 61.70674100000001
$ws.Cells.Item(12, 9).Value = 0.7535372480335895
$ws.Cells.Item(12, 10).Value = 0.7535372480335896
$ws.Cells.Item(12, 13).Value = 15.42507366666667
$ws.Cells.Item(12, 14).Value = 46.275221
$ws.Cells.Item(12, 15).Value = 0.1737373292074942
$ws.Cells.Item(12, 16).Value = 0.1737373292074942
$ws.Cells.Item(12, 17).Value = 317.2770085516402
$ws.Cells.Item(12, 18).Value = 2855.493076964762
$ws.Cells.Item(12, 19).Value = 0.1309175489317209
$ws.Cells.Item(12, 20).Value = 0.1309175489317209

$ws.Cells.Item(13, 7).Value = 20.56891366666667
$ws.Cells.Item(13, 8).Value = 61.70674100000001
$ws.Cells.Item(13, 9).Value = 0.7535372480335895
$ws.Cells.Item(13, 10).Value = 0.7535372480335896
$ws.Cells.Item(13, 13).Value = 13.804248
$ws.Cells.Item(13, 14).Value = 41.412744
$ws.Cells.Item(13, 15).Value = 0.1554814732859661
$ws.Cells.Item(13, 16).Value = 0.1554814732859661
$ws.Cells.Item(13, 17).Value = 283.9383853452561
$ws.Cells.Item(13, 18).Value = 2555.445468107304
$ws.Cells.Item(13, 19).Value = 0.117161081500115
$ws.Cells.Item(13, 20).Value = 0.117161081500115

$ws.Cells.Item(14, 7).Value = 1.155605333333333
$ws.Cells.Item(14, 8).Value = 3.466816
$ws.Cells.Item(14, 9).Value = 0.04233532586138062
$ws.Cells.Item(14, 10).Value = 0.04233532586138063
$ws.Cells.Item(14, 13).Value = 15.090721
$ws.Cells.Item(14, 14).Value = 45.272163
$ws.Cells.Item(14, 15).Value = 0.169971412714946
$ws.Cells.Item(14, 16).Value = 0.169971412714946
$ws.Cells.Item(14, 17).Value = 17.43891767144533
$ws.Cells.Item(14, 18).Value = 156.950259043008
$ws.Cells.Item(14, 19).Value = 0.007195795144406453
$ws.Cells.Item(14, 20).Value = 0.007195795144406453

$ws.Cells.Item(15, 7).Value = 1.155605333333333
$ws.Cells.Item(15, 8).Value = 3.466816
$ws.Cells.Item(15, 9).Value = 0.04233532586138062
$ws.Cells.Item(15, 10).Value = 0.04233532586138063
$ws.Cells.Item(15, 14).Value = 60.45961299999999
$ws.Cells.Item(15, 15).Value = 0.226991713071207
$ws.Cells.Item(15, 16).Value = 0.226991713071207
$ws.Cells.Item(15, 17).Value = 23.28915041135644
$ws.Cells.Item(15, 18).Value = 209.6023537022079
$ws.Cells.Item(15, 19).Value = 0.009609768140702558
$ws.Cells.Item(15, 20).Value = 0.00960976814070256

$ws.Cells.Item(16, 7).Value = 1.155605333333333
$ws.Cells.Item(16, 8).Value = 3.466816
$ws.Cells.Item(16, 9).Value = 0.04233532586138062
$ws.Cells.Item(16, 10).Value = 0.04233532586138063
$ws.Cells.Item(16, 13).Value = 18.66868666666667
$ws.Cells.Item(16, 14).Value = 56.00606
$ws.Cells.Item(16, 15).Value = 0.2102711359030499
$ws.Cells.Item(16, 16).Value = 0.2102711359030499
$ws.Cells.Item(16, 17).Value = 21.57363387832888
$ws.Cells.Item(16, 18).Value = 194.16270490496
$ws.Cells.Item(16, 19).Value = 0.008901897057698267
$ws.Cells.Item(16, 20).Value = 0.00890189705769827

$ws.Cells.Item(17, 7).Value = 1.155605333333333
$ws.Cells.Item(17, 8).Value = 3.466816
$ws.Cells.Item(17, 9).Value = 0.04233532586138062
$ws.Cells.Item(17, 10).Value = 0.04233532586138063
$ws.Cells.Item(17, 13).Value = 5.641943333333334
$ws.Cells.Item(17, 14).Value = 16.92583
$ws.Cells.Item(17, 15).Value = 0.06354693581733691
$ws.Cells.Item(17, 16).Value = 0.06354693581733692
$ws.Cells.Item(17, 17).Value = 6.519859806364444
$ws.Cells.Item(17, 18).Value = 58.67873825728
$ws.Cells.Item(17, 19).Value = 0.002690280235319197
$ws.Cells.Item(17, 20).Value = 0.002690280235319198

$ws.Cells.Item(18, 7).Value = 1.155605333333333
$ws.Cells.Item(18, 8).Value = 3.466816
$ws.Cells.Item(18, 9).Value = 0.04233532586138062
$ws.Cells.Item(18, 10).Value = 0.04233532586138063
$ws.Cells.Item(18, 13).Value = 15.42507366666667
$ws.Cells.Item(18, 14).Value = 46.275221
$ws.Cells.Item(18, 15).Value = 0.1737373292074942
$ws.Cells.Item(18, 16).Value = 0.1737373292074942
$ws.Cells.Item(18, 17).Value = 17.82529739625955
$ws.Cells.Item(18, 18).Value = 160.427676566336
$ws.Cells.Item(18, 19).Value = 0.007355226446285226
$ws.Cells.Item(18, 20).Value = 0.007355226446285228

$ws.Cells.Item(19, 7).Value = 1.155605333333333
$ws.Cells.Item(19, 8).Value = 3.466816
$ws.Cells.Item(19, 9).Value = 0.04233532586138062
$ws.Cells.Item(19, 10).Value = 0.04233532586138063
$ws.Cells.Item(19, 13).Value = 13.804248
$ws.Cells.Item(19, 14).Value = 41.412744
$ws.Cells.Item(19, 15).Value = 0.1554814732859661
$ws.Cells.Item(19, 16).Value = 0.1554814732859661
$ws.Cells.Item(19, 17).Value = 15.952262611456
$ws.Cells.Item(19, 18).Value = 143.570363503104
$ws.Cells.Item(19, 19).Value = 0.006582358836968921
$ws.Cells.Item(19, 20).Value = 0.006582358836968923
